$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.472.01"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "'2.095.10"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D5").Value = "'330.14"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.5213"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "'0.4361"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").Value = "'53.91"
$ws.Range("E9").Value = "  +14.91%  "
$ws.Range("D10").Value = "'0.08867"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("D11").Value = "'1.155"
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").Value = "'24.40"
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("D13").Value = "'2.092.89"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "'6.673"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").Value = "'7.668"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "'95.82"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'0.00001122"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "'0.06587"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "'19.28"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'6.248"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("D23").Value = "'30.496.37"
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("D24").Value = "'12.25"
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("D25").Value = "'2.341"
$ws.Range("E25").Value = "  +3.76%  "
$ws.Range("D26").Value = "'2.336.84"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").Value = "'22.26"
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("D28").Value = "'2.564"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "'162.55"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "'131.77"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").Value = "'1.182"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "'1.645"
$ws.Range("E33").Value = "  +7.21%  "
$ws.Range("D34").Value = "'6.135"
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").Value = "'3.894"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").Value = "'10.00"
$ws.Range("E36").Value = "  +4.96%  "
$ws.Range("D37").Value = "'0.02572"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").Value = "'0.06809"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").Value = "'12.75"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").Value = "'5.441"
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("D41").Value = "'0.2256"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "'0.6874"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "'1.258"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "'0.6361"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").Value = "'13.87"
$ws.Range("E46").Value = "  -3.16%  "
$ws.Range("D47").Value = "'2.195"
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("D48").Value = "'3.627"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "'1.237"
$ws.Range("E49").Value = "  +8.20%  "
$ws.Range("D50").Value = "'1.240"
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("D51").Value = "'81.70"
$ws.Range("E51").Value = "  -1.74%  "
